$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3810755   # was 4250390.5
$ws.Range("J17").Value = 4250400   # was 4420386
$ws.Range("L17").Value = 12751200   # was 13261158
$ws.Range("N17").Value = -12751536   # was -13261494
$ws.Range("H68").Value = 18750   # was 30000
$ws.Range("J68").Value = 18750   # was 30000
$ws.Range("L68").Value = 18750   # was 30000
$ws.Range("N68").Value = -20248   # was -31498
$ws.Range("H71").Value = 18750   # was 30000
$ws.Range("J71").Value = 18750   # was 30000
$ws.Range("L71").Value = 56250   # was 90000
$ws.Range("N71").Value = -63738   # was -97488
$ws.Range("H106").Value = 7854   # was 7633.1055
$ws.Range("I106").Value = 8302.143   # was 7768.6665
$ws.Range("J106").Value = 6599.2   # was 7124.75
$ws.Range("K106").Value = 8302.143   # was 7768.6665
$ws.Range("L106").Value = 6599.2   # was 7124.75
$ws.Range("M106").Value = -7671.143   # was -7137.6665
$ws.Range("N106").Value = -7861.2   # was -8386.75
$ws.Range("H107").Value = 7085.625   # was 6283.8335
$ws.Range("I107").Value = 18504.834   # was 12329.777
$ws.Range("J107").Value = 234.1   # was 237.88889
$ws.Range("K107").Value = 18504.834   # was 12329.777
$ws.Range("L107").Value = 234.1   # was 237.88889
$ws.Range("M107").Value = -16584.834   # was -10409.777
$ws.Range("N107").Value = -4074.1   # was -4077.88889
$ws.Range("H113").Value = 2659.889   # was 2670.4443
$ws.Range("I113").Value = 2229.2307   # was 2243.8462
$ws.Range("K113").Value = 2229.2307   # was 2243.8462
$ws.Range("M113").Value = 1024.7693   # was 1010.1538
$ws.Range("H138").Value = 4568735.5   # was 4568708
$ws.Range("I138").Value = 9805652   # was 9525489
$ws.Range("J138").Value = 3218.5386   # was 3251.658
$ws.Range("K138").Value = 29416956   # was 28576467
$ws.Range("L138").Value = 9655.6158   # was 9754.974
$ws.Range("M138").Value = -29411816   # was -28571327
$ws.Range("N138").Value = -19935.6158   # was -20034.974

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2118.0635   # was 2391.8867
$ws.Range("I61").Value = 2034.8518   # was 2311.422
$ws.Range("J61").Value = 2617.3333   # was 2844.5
$ws.Range("K61").Value = 2034.8518   # was 2311.422
$ws.Range("L61").Value = 2617.3333   # was 2844.5
$ws.Range("M61").Value = -1822.8518   # was -2099.422
$ws.Range("N61").Value = -3041.3333   # was -3268.5
$ws.Range("H97").Value = 4950.5454   # was 8260.846
$ws.Range("I97").Value = 5772.778   # was 10228
$ws.Range("J97").Value = 1250.5   # was 1703.6666
$ws.Range("K97").Value = 5772.778   # was 10228
$ws.Range("L97").Value = 1250.5   # was 1703.6666
$ws.Range("M97").Value = -5276.778   # was -9732
$ws.Range("N97").Value = -2242.5   # was -2695.6666
$ws.Range("H110").Value = 1320.3334   # was 789.64703
$ws.Range("I110").Value = 1320.3334   # was 534.38464
$ws.Range("J110").Value = 0   # was 1619.25
$ws.Range("K110").Value = 1320.3334   # was 534.38464
$ws.Range("L110").Value = 0   # was 1619.25
$ws.Range("M110").Value = 724.6666   # was 1510.61536
$ws.Range("N110").Value = ""   # was -5709.25
$ws.Range("H136").Value = 2118.0635   # was 2391.8867
$ws.Range("I136").Value = 2034.8518   # was 2311.422
$ws.Range("J136").Value = 2617.3333   # was 2844.5
$ws.Range("K136").Value = 6104.555399999999   # was 6934.266
$ws.Range("L136").Value = 7851.999899999999   # was 8533.5
$ws.Range("M136").Value = -3554.555399999999   # was -4384.266
$ws.Range("N136").Value = -12951.9999   # was -13633.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1251.2858   # was 711.875
$ws.Range("I99").Value = 937.5   # was 660.7692
$ws.Range("J99").Value = 1669.6666   # was 933.3333
$ws.Range("K99").Value = 937.5   # was 660.7692
$ws.Range("L99").Value = 1669.6666   # was 933.3333
$ws.Range("M99").Value = 560.5   # was 837.2308
$ws.Range("N99").Value = -4665.6666   # was -3929.3333
$ws.Range("H105").Value = 3427.923   # was 4047.8484
$ws.Range("I105").Value = 1852.4375   # was 2389.875
$ws.Range("J105").Value = 4523.913   # was 4578.4
$ws.Range("K105").Value = 1852.4375   # was 2389.875
$ws.Range("L105").Value = 4523.913   # was 4578.4
$ws.Range("M105").Value = -105.4375   # was -642.875
$ws.Range("N105").Value = -8017.913   # was -8072.4
$ws.Range("H134").Value = 3319.3274   # was 3617.102
$ws.Range("I134").Value = 2004.6552   # was 2173
$ws.Range("J134").Value = 4785.6924   # was 5121.375
$ws.Range("K134").Value = 6013.9656   # was 6519
$ws.Range("L134").Value = 14357.0772   # was 15364.125
$ws.Range("M134").Value = -3478.9656   # was -3984
$ws.Range("N134").Value = -19427.0772   # was -20434.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1704.3334   # was 1756.7646
$ws.Range("I99").Value = 1370.9286   # was 1489.1666
$ws.Range("J99").Value = 2371.1428   # was 2399
$ws.Range("K99").Value = 1370.9286   # was 1489.1666
$ws.Range("L99").Value = 2371.1428   # was 2399
$ws.Range("M99").Value = 127.0714   # was 8.833399999999983
$ws.Range("N99").Value = -5367.1428   # was -5395
$ws.Range("H126").Value = 1704.3334   # was 1756.7646
$ws.Range("I126").Value = 1370.9286   # was 1489.1666
$ws.Range("J126").Value = 2371.1428   # was 2399
$ws.Range("K126").Value = 4112.7858   # was 4467.4998
$ws.Range("L126").Value = 7113.428400000001   # was 7197
$ws.Range("M126").Value = -1642.7858   # was -1997.4998
$ws.Range("N126").Value = -12053.4284   # was -12137

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 86.72727   # was 117.14286
$ws.Range("I8").Value = 86.72727   # was 117.14286
$ws.Range("K8").Value = 260.18181   # was 351.42858
$ws.Range("M8").Value = -121.18181   # was -212.42858
$ws.Range("H87").Value = 10942.842   # was 13653.846
$ws.Range("I87").Value = 6762.615   # was 9687.5
$ws.Range("K87").Value = 20287.845   # was 29062.5
$ws.Range("M87").Value = -19039.845   # was -27814.5
$ws.Range("H90").Value = 10942.842   # was 13653.846
$ws.Range("I90").Value = 6762.615   # was 9687.5
$ws.Range("K90").Value = 60863.535   # was 87187.5
$ws.Range("M90").Value = -54623.535   # was -80947.5
$ws.Range("H121").Value = 350   # was 464.16666
$ws.Range("I121").Value = 350   # was 437
$ws.Range("J121").Value = 0   # was 600
$ws.Range("K121").Value = 1050   # was 1311
$ws.Range("L121").Value = 0   # was 1800
$ws.Range("M121").Value = 260   # was -1
$ws.Range("N121").Value = ""   # was -4420
$ws.Range("H132").Value = 1410.3529   # was 1525.5625
$ws.Range("I132").Value = 670.7273   # was 684.4
$ws.Range("J132").Value = 2766.3333   # was 2927.5
$ws.Range("K132").Value = 6036.545700000001   # was 6159.599999999999
$ws.Range("L132").Value = 24896.9997   # was 26347.5
$ws.Range("M132").Value = -3506.545700000001   # was -3629.599999999999
$ws.Range("N132").Value = -29956.9997   # was -31407.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 201229.8   # was 201012.8
$ws.Range("I113").Value = 334183.66   # was 334117
$ws.Range("J113").Value = 1799   # was 1356.5
$ws.Range("K113").Value = 334183.66   # was 334117
$ws.Range("L113").Value = 1799   # was 1356.5
$ws.Range("M113").Value = -332013.66   # was -331947
$ws.Range("N113").Value = -6139   # was -5696.5
$ws.Range("H122").Value = 1901.2222   # was 2115.36
$ws.Range("I122").Value = 1855.5625   # was 1986
$ws.Range("J122").Value = 1967.6364   # was 2309.4
$ws.Range("K122").Value = 5566.6875   # was 5958
$ws.Range("L122").Value = 5902.9092   # was 6928.200000000001
$ws.Range("M122").Value = -3116.6875   # was -3508
$ws.Range("N122").Value = -10802.9092   # was -11828.2
$ws.Range("H132").Value = 5227.9688   # was 5374.032
$ws.Range("I132").Value = 5555   # was 5741.731
$ws.Range("K132").Value = 16665   # was 17225.193
$ws.Range("M132").Value = -14135   # was -14695.193

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5291.3477   # was 6163
$ws.Range("I40").Value = 5708.25   # was 8913.571
$ws.Range("J40").Value = 4836.5454   # was 4558.5
$ws.Range("K40").Value = 5708.25   # was 8913.571
$ws.Range("L40").Value = 4836.5454   # was 4558.5
$ws.Range("M40").Value = -5572.25   # was -8777.571
$ws.Range("N40").Value = -5108.5454   # was -4830.5
$ws.Range("H68").Value = 1506.1666   # was 1133.6666
$ws.Range("I68").Value = 1425.5   # was 900
$ws.Range("J68").Value = 2151.5   # was 1200.4286
$ws.Range("K68").Value = 1425.5   # was 900
$ws.Range("L68").Value = 2151.5   # was 1200.4286
$ws.Range("M68").Value = -676.5   # was -151
$ws.Range("N68").Value = -3649.5   # was -2698.4286
$ws.Range("H71").Value = 1506.1666   # was 1133.6666
$ws.Range("I71").Value = 1425.5   # was 900
$ws.Range("J71").Value = 2151.5   # was 1200.4286
$ws.Range("K71").Value = 7127.5   # was 4500
$ws.Range("L71").Value = 10757.5   # was 6002.143
$ws.Range("M71").Value = -3383.5   # was -756
$ws.Range("N71").Value = -18245.5   # was -13490.143
$ws.Range("H82").Value = 1304.3043   # was 1248.4762
$ws.Range("I82").Value = 1150.4445   # was 1132.5264
$ws.Range("J82").Value = 1858.2   # was 2350
$ws.Range("K82").Value = 1150.4445   # was 1132.5264
$ws.Range("L82").Value = 1858.2   # was 2350
$ws.Range("M82").Value = -789.4445000000001   # was -771.5264
$ws.Range("N82").Value = -2580.2   # was -3072
$ws.Range("H85").Value = 1304.3043   # was 1248.4762
$ws.Range("I85").Value = 1150.4445   # was 1132.5264
$ws.Range("J85").Value = 1858.2   # was 2350
$ws.Range("K85").Value = 1150.4445   # was 1132.5264
$ws.Range("L85").Value = 1858.2   # was 2350
$ws.Range("M85").Value = 97.55549999999994   # was 115.4736
$ws.Range("N85").Value = -4354.2   # was -4846
$ws.Range("H122").Value = 7106.207   # was 7276.357
$ws.Range("I122").Value = 8640.25   # was 8722.5
$ws.Range("J122").Value = 6521.8096   # was 6697.9
$ws.Range("K122").Value = 25920.75   # was 26167.5
$ws.Range("L122").Value = 19565.4288   # was 20093.7
$ws.Range("M122").Value = -23470.75   # was -23717.5
$ws.Range("N122").Value = -24465.4288   # was -24993.7
$ws.Range("H132").Value = 8778661   # was 15635466
$ws.Range("I132").Value = 4869.2354   # was 9899.538
$ws.Range("J132").Value = 21748612   # was 26326642
$ws.Range("K132").Value = 14607.7062   # was 29698.614
$ws.Range("L132").Value = 65245836   # was 78979926
$ws.Range("M132").Value = -12077.7062   # was -27168.614
$ws.Range("N132").Value = -65250896   # was -78984986

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5406.6665   # was 5254.5454
$ws.Range("I62").Value = 5280   # was 5685.7144
$ws.Range("J62").Value = 5660   # was 4500
$ws.Range("K62").Value = 5280   # was 5685.7144
$ws.Range("L62").Value = 5660   # was 4500
$ws.Range("M62").Value = -4656   # was -5061.7144
$ws.Range("N62").Value = -6908   # was -5748
$ws.Range("H65").Value = 5406.6665   # was 5254.5454
$ws.Range("I65").Value = 5280   # was 5685.7144
$ws.Range("J65").Value = 5660   # was 4500
$ws.Range("K65").Value = 26400   # was 28428.572
$ws.Range("L65").Value = 28300   # was 22500
$ws.Range("M65").Value = -23280   # was -25308.572
$ws.Range("N65").Value = -34540   # was -28740
$ws.Range("H136").Value = 1446.6305   # was 2085.8667
$ws.Range("I136").Value = 874.21875   # was 1331.1052
$ws.Range("J136").Value = 2755   # was 3389.5454
$ws.Range("K136").Value = 2622.65625   # was 3993.3156
$ws.Range("L136").Value = 8265   # was 10168.6362
$ws.Range("M136").Value = -72.65625   # was -1443.3156
$ws.Range("N136").Value = -13365   # was -15268.6362
